# Saldo_guide.xlsx update:
#   - Reference date moved forward one day (2024-12-09 -> 2024-12-10),
#     reflected both in the worksheet's generated name and in every
#     "Dt. Referencia" (column G) value on the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Dt. Referencia") holds serial date 45635 (2024-12-09) for
# every data row (2 through 274); bump each to 45636 (2024-12-10).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45636
}

# Sheet name encodes the same export timestamp; rename it to match the
# new reference date/time.
$ws.Name = "IClientBalance-20241210-075931-"
